$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-10-21 00:36:06"
$wsZhCn.Range("P2").Value = ""

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-10-21 00:36:24"
$wsDeDe.Range("P2").Value = ""
